$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.441378666666667
$ws.Range("H2").Value = 7.324135999999999
$ws.Range("I2").Value = 0.1119936059016048
$ws.Range("J2").Value = 0.1119936059016048
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.07474266666667
$ws.Range("N2").Value = 36.224228
$ws.Range("O2").Value = 0.08154942646895191
$ws.Range("P2").Value = 0.08154942646895191
$ws.Range("Q2").Value = 29.47901915188977
$ws.Range("R2").Value = 265.3111723670079
$ws.Range("S2").Value = 0.0091330143294657
$ws.Range("T2").Value = 0.009133014329465698

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.441378666666667
$ws.Range("H3").Value = 7.324135999999999
$ws.Range("I3").Value = 0.1119936059016048
$ws.Range("J3").Value = 0.1119936059016048
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 16.614382
$ws.Range("N3").Value = 49.84314599999999
$ws.Range("O3").Value = 0.1122088777077108
$ws.Range("P3").Value = 0.1122088777077108
$ws.Range("Q3").Value = 40.56199777465066
$ws.Range("R3").Value = 365.0579799718559
$ws.Range("S3").Value = 0.01256667682865873
$ws.Range("T3").Value = 0.01256667682865873

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.441378666666667
$ws.Range("H4").Value = 7.324135999999999
$ws.Range("I4").Value = 0.1119936059016048
$ws.Range("J4").Value = 0.1119936059016048
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 43.86740866666667
$ws.Range("N4").Value = 131.602226
$ws.Range("O4").Value = 0.2962681786437903
$ws.Range("P4").Value = 0.2962681786437903
$ws.Range("Q4").Value = 107.0969556807485
$ws.Range("R4").Value = 963.872601126736
$ws.Range("S4").Value = 0.0331801416402189
$ws.Range("T4").Value = 0.0331801416402189

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.441378666666667
$ws.Range("H5").Value = 7.324135999999999
$ws.Range("I5").Value = 0.1119936059016048
$ws.Range("J5").Value = 0.1119936059016048
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 69.577158
$ws.Range("N5").Value = 208.731474
$ws.Range("O5").Value = 0.4699046171727648
$ws.Range("P5").Value = 0.4699046171727648
$ws.Range("Q5").Value = 169.864189228496
$ws.Range("R5").Value = 1528.777703056464
$ws.Range("S5").Value = 0.0526263125069911
$ws.Range("T5").Value = 0.0526263125069911

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.441378666666667
$ws.Range("H6").Value = 7.324135999999999
$ws.Range("I6").Value = 0.1119936059016048
$ws.Range("J6").Value = 0.1119936059016048
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.932863999999999
$ws.Range("N6").Value = 17.798592
$ws.Range("O6").Value = 0.04006890000678209
$ws.Range("P6").Value = 0.04006890000678209
$ws.Range("Q6").Value = 14.48436760183466
$ws.Range("R6").Value = 130.359308416512
$ws.Range("S6").Value = 0.004487460596270362
$ws.Range("T6").Value = 0.004487460596270362

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.88577866666667
$ws.Range("H7").Value = 38.657336
$ws.Range("I7").Value = 0.5911106037886134
$ws.Range("J7").Value = 0.5911106037886134
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.07474266666667
$ws.Range("N7").Value = 36.224228
$ws.Range("O7").Value = 0.08154942646895191
$ws.Range("P7").Value = 0.08154942646895191
$ws.Range("Q7").Value = 155.5924614596231
$ws.Range("R7").Value = 1400.332153136608
$ws.Range("S7").Value = 0.04820473071867729
$ws.Range("T7").Value = 0.04820473071867729

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.88577866666667
$ws.Range("H8").Value = 38.657336
$ws.Range("I8").Value = 0.5911106037886134
$ws.Range("J8").Value = 0.5911106037886134
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.614382
$ws.Range("N8").Value = 49.84314599999999
$ws.Range("O8").Value = 0.1122088777077108
$ws.Range("P8").Value = 0.1122088777077108
$ws.Range("Q8").Value = 214.0892491354506
$ws.Range("R8").Value = 1926.803242219056
$ws.Range("S8").Value = 0.06632785745224762
$ws.Range("T8").Value = 0.06632785745224762

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.88577866666667
$ws.Range("H9").Value = 38.657336
$ws.Range("I9").Value = 0.5911106037886134
$ws.Range("J9").Value = 0.5911106037886134
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 43.86740866666667
$ws.Range("N9").Value = 131.602226
$ws.Range("O9").Value = 0.2962681786437903
$ws.Range("P9").Value = 0.2962681786437903
$ws.Range("Q9").Value = 565.2657187588818
$ws.Range("R9").Value = 5087.391468829936
$ws.Range("S9").Value = 0.1751272619614837
$ws.Range("T9").Value = 0.1751272619614837

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.88577866666667
$ws.Range("H10").Value = 38.657336
$ws.Range("I10").Value = 0.5911106037886134
$ws.Range("J10").Value = 0.5911106037886134
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 69.577158
$ws.Range("N10").Value = 208.731474
$ws.Range("O10").Value = 0.4699046171727648
$ws.Range("P10").Value = 0.4699046171727648
$ws.Range("Q10").Value = 896.555858243696
$ws.Range("R10").Value = 8069.002724193264
$ws.Range("S10").Value = 0.2777656019800502
$ws.Range("T10").Value = 0.2777656019800502

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.88577866666667
$ws.Range("H11").Value = 38.657336
$ws.Range("I11").Value = 0.5911106037886134
$ws.Range("J11").Value = 0.5911106037886134
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.932863999999999
$ws.Range("N11").Value = 17.798592
$ws.Range("O11").Value = 0.04006890000678209
$ws.Range("P11").Value = 0.04006890000678209
$ws.Range("Q11").Value = 76.44957236343465
$ws.Range("R11").Value = 688.046151270912
$ws.Range("S11").Value = 0.02368515167615453
$ws.Range("T11").Value = 0.02368515167615453

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.472111000000001
$ws.Range("H12").Value = 19.416333
$ws.Range("I12").Value = 0.2968957903097819
$ws.Range("J12").Value = 0.2968957903097818
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.07474266666667
$ws.Range("N12").Value = 36.224228
$ws.Range("O12").Value = 0.08154942646895191
$ws.Range("P12").Value = 0.08154942646895191
$ws.Range("Q12").Value = 78.14907483510267
$ws.Range("R12").Value = 703.341673515924
$ws.Range("S12").Value = 0.02421168142080892
$ws.Range("T12").Value = 0.02421168142080892

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.472111000000001
$ws.Range("H13").Value = 19.416333
$ws.Range("I13").Value = 0.2968957903097819
$ws.Range("J13").Value = 0.2968957903097818
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.614382
$ws.Range("N13").Value = 49.84314599999999
$ws.Range("O13").Value = 0.1122088777077108
$ws.Range("P13").Value = 0.1122088777077108
$ws.Range("Q13").Value = 107.530124500402
$ws.Range("R13").Value = 967.7711205036179
$ws.Range("S13").Value = 0.03331434342680447
$ws.Range("T13").Value = 0.03331434342680446

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.472111000000001
$ws.Range("H14").Value = 19.416333
$ws.Range("I14").Value = 0.2968957903097819
$ws.Range("J14").Value = 0.2968957903097818
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 43.86740866666667
$ws.Range("N14").Value = 131.602226
$ws.Range("O14").Value = 0.2962681786437903
$ws.Range("P14").Value = 0.2962681786437903
$ws.Range("Q14").Value = 283.9147381730287
$ws.Range("R14").Value = 2555.232643557258
$ws.Range("S14").Value = 0.08796077504208777
$ws.Range("T14").Value = 0.08796077504208775

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.472111000000001
$ws.Range("H15").Value = 19.416333
$ws.Range("I15").Value = 0.2968957903097819
$ws.Range("J15").Value = 0.2968957903097818
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 69.577158
$ws.Range("N15").Value = 208.731474
$ws.Range("O15").Value = 0.4699046171727648
$ws.Range("P15").Value = 0.4699046171727648
$ws.Range("Q15").Value = 450.3110896405381
$ws.Range("R15").Value = 4052.799806764842
$ws.Range("S15").Value = 0.1395127026857235
$ws.Range("T15").Value = 0.1395127026857235

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.472111000000001
$ws.Range("H16").Value = 19.416333
$ws.Range("I16").Value = 0.2968957903097819
$ws.Range("J16").Value = 0.2968957903097818
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.932863999999999
$ws.Range("N16").Value = 17.798592
$ws.Range("O16").Value = 0.04006890000678209
$ws.Range("P16").Value = 0.04006890000678209
$ws.Range("Q16").Value = 38.398154355904
$ws.Range("R16").Value = 345.583389203136
$ws.Range("S16").Value = 0.01189628773435719
$ws.Range("T16").Value = 0.01189628773435719
